$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text values (no risk of Excel reinterpreting them as a number/date),
# safe to assign directly.
$ws.Range("E2").Value = "SHUBR000"
$ws.Range("G2").Value = "SHUMAN"
$ws.Range("H2").Value = "BRIAN"
$ws.Range("I2").Value = "J"
$ws.Range("L2").Value = "PO BOX 1174"
$ws.Range("M2").Value = "LA"
$ws.Range("N2").Value = "ABITA SPRINGS"
$ws.Range("R2").Value = "Stratton Beatrous Grisoli M.D.,"
$ws.Range("S2").Value = "Baldone Reina Dermatology, APMC"
$ws.Range("T2").Value = "BRIAN SHUMAN"
$ws.Range("W2").Value = "Medicare of Louisiana"
$ws.Range("Y2").Value = "4NT5FH0HP64"
$ws.Range("AD2").Value = "BRIAN SHUMAN"
$ws.Range("AE2").Value = "Self"
$ws.Range("AG2").Value = "Physicians Mutual Insurance Company"
$ws.Range("AI2").Value = "H730052094"

# Values that look numeric/date-like must stay stored as text (matching the
# source data, e.g. requisition numbers with leading zeros or a trailing
# dot). Enter them as a literal-text formula, then paste the result back as
# a static value so the cell ends up a plain text (shared-string) cell with
# no formula and no number-format style attached.
$textLikeNumbers = @{
    "C2"  = "33556"
    "D2"  = "32847311."
    "J2"  = "1956-12-31"
    "O2"  = "704201174"
    "P2"  = "9856309533"
    "V2"  = "1956-12-31"
    "AF2" = "1956-12-31"
}
foreach ($addr in $textLikeNumbers.Keys) {
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $textLikeNumbers[$addr] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# The old primary-insurance street address / city / state / zip no longer
# apply to the new record, so clear those cells entirely.
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").ClearContents()
